$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.959.67'
$ws.Range("E2").Value = '  -0.91%  '

$ws.Range("D3").Value = '3.407.94'
$ws.Range("E3").Value = '  -0.74%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '409.48'
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("D6").Value = '128.88'
$ws.Range("E6").Value = '  -1.56%  '

$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +6.38%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.734'
$ws.Range("E9").Value = '  +5.23%  '

$ws.Range("E10").Value = '  -0.40%  '

$ws.Range("D11").Value = '43.03'
$ws.Range("E11").Value = '  +1.97%  '

$ws.Range("D12").Value = '0.0000219'
$ws.Range("E12").Value = '  +31.22%  '

$ws.Range("D13").Value = '9.30'
$ws.Range("E13").Value = '  +9.56%  '

$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("D15").Value = '21.37'
$ws.Range("E15").Value = '  +7.31%  '

$ws.Range("D16").Value = '3.947.52'
$ws.Range("E16").Value = '  -0.83%  '

$ws.Range("D17").Value = '3.414.49'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("D18").Value = '12.54'
$ws.Range("E18").Value = '  +8.49%  '

$ws.Range("E19").Value = '  +6.67%  '

$ws.Range("D20").Value = '61.982.45'
$ws.Range("E20").Value = '  -0.95%  '

$ws.Range("D21").Value = '447.61'
$ws.Range("E21").Value = '  +42.20%  '

$ws.Range("D22").Value = '91.98'
$ws.Range("E22").Value = '  +8.78%  '

$ws.Range("E23").Value = '  +0.21%  '

$ws.Range("D24").Value = '13.23'
$ws.Range("E24").Value = '  +2.82%  '

$ws.Range("D25").Value = '3.29'
$ws.Range("E25").Value = '  +3.41%  '

$ws.Range("D26").Value = '9.39'
$ws.Range("E26").Value = '  +14.79%  '

$ws.Range("D27").Value = '33.13'
$ws.Range("E27").Value = '  +11.03%  '

$ws.Range("D28").Value = '4.80'
$ws.Range("E28").Value = '  +1.24%  '

$ws.Range("D29").Value = '7.73'
$ws.Range("E29").Value = '  -1.00%  '

$ws.Range("D30").Value = '2.70'
$ws.Range("E30").Value = '  -1.06%  '

$ws.Range("D31").Value = '12.00'
$ws.Range("E31").Value = '  +4.75%  '

$ws.Range("E32").Value = '  -1.67%  '

$ws.Range("E33").Value = '  -0.90%  '

$ws.Range("D34").Value = '42.71'
$ws.Range("E34").Value = '  -4.02%  '

$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").Value = '0.0505'
$ws.Range("E36").Value = '  +3.94%  '

$ws.Range("D37").Value = '53.88'
$ws.Range("E37").Value = '  +4.00%  '

$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("D39").Value = '3.39'
$ws.Range("E39").Value = '  +1.48%  '

$ws.Range("D40").Value = '0.136'
$ws.Range("E40").Value = '  +7.89%  '

$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("D42").Value = '0.318'
$ws.Range("E42").Value = '  -1.69%  '

$ws.Range("D43").Value = '143.96'
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("D44").Value = '4.32'
$ws.Range("E44").Value = '  +9.71%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = '2.00'
$ws.Range("E45").Value = '  +0.77%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '2.55'
$ws.Range("E46").Value = '  +14.44%  '

$ws.Range("D47").Value = '16.64'
$ws.Range("E47").Value = '  -1.91%  '

$ws.Range("D48").Value = '0.152'
$ws.Range("E48").Value = '  +24.71%  '

$ws.Range("D49").Value = '22.50'
$ws.Range("E49").Value = '  +5.75%  '

$ws.Range("D50").Value = '2.17'
$ws.Range("E50").Value = '  +8.25%  '

$ws.Range("D51").Value = '3.750.24'
$ws.Range("E51").Value = '  -0.66%  '
